$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 69 (shifts the existing rows 69..88 down to 70..89)
$ws.Rows.Item(69).Insert()

# Populate the newly inserted row 69 with the new weekly price record
$ws.Cells.Item(69, 1).Value  = 10
$ws.Cells.Item(69, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(69, 3).Value  = "La Araucanía"
$ws.Cells.Item(69, 4).Value  = 44736
$ws.Cells.Item(69, 5).Value  = 9
$ws.Cells.Item(69, 6).Value  = 100112035
$ws.Cells.Item(69, 7).Value  = "Bruselas (repollito)"
$ws.Cells.Item(69, 8).Value  = "Sin especificar"
$ws.Cells.Item(69, 9).Value  = "Primera"
$ws.Cells.Item(69, 10).Value = 30
$ws.Cells.Item(69, 11).Value = 26000
$ws.Cells.Item(69, 12).Value = 26000
$ws.Cells.Item(69, 13).Value = 26000
$ws.Cells.Item(69, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(69, 15).Value = "Región Metropolitana"
$ws.Cells.Item(69, 16).Value = 2600
$ws.Cells.Item(69, 17).Value = 10
$ws.Cells.Item(69, 18).Value = "Hortaliza"
